# CPL_JLCPCB_ni_arena_12-12_v0p1_r1.xlsx
# Add new placement rows for J1, P1-P12, SW1 connectors/pogo-pins and the
# reset switch, replacing the trailing block of empty rows that padded the
# sheet out to the max row (1048576).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 14 trailing blank rows (1048563:1048576) that existed only to
# carry row-height formatting; the dimension will shrink back down once
# they're gone.
$ws.Rows("1048563:1048576").Delete()

# Stamp the formatting (style) of an existing data row onto the new block
# of rows so every new cell picks up the same cell style (s="2") as the
# rest of the CPL table, then fill in the values.
$ws.Range("A2:E2").Copy()
$ws.Range("A70:E83").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(70, 1).Value = "J1"
$ws.Cells.Item(70, 2).Value = 304.6
$ws.Cells.Item(70, 3).Value = -116.459
$ws.Cells.Item(70, 4).Value = "top"
$ws.Cells.Item(70, 5).Value = -90

$ws.Cells.Item(71, 1).Value = "P1"
$ws.Cells.Item(71, 2).Value = 265.474
$ws.Cells.Item(71, 3).Value = -177.8
$ws.Cells.Item(71, 4).Value = "top"
$ws.Cells.Item(71, 5).Value = 90

$ws.Cells.Item(72, 1).Value = "P2"
$ws.Cells.Item(72, 2).Value = 253.7279
$ws.Cells.Item(72, 3).Value = -133.963
$ws.Cells.Item(72, 4).Value = "top"
$ws.Cells.Item(72, 5).Value = 120

$ws.Cells.Item(73, 1).Value = "P3"
$ws.Cells.Item(73, 2).Value = 221.637
$ws.Cells.Item(73, 3).Value = -101.8721
$ws.Cells.Item(73, 4).Value = "top"
$ws.Cells.Item(73, 5).Value = 150

$ws.Cells.Item(74, 1).Value = "P4"
$ws.Cells.Item(74, 2).Value = 177.8
$ws.Cells.Item(74, 3).Value = -90.126
$ws.Cells.Item(74, 4).Value = "top"
$ws.Cells.Item(74, 5).Value = 180

$ws.Cells.Item(75, 1).Value = "P5"
$ws.Cells.Item(75, 2).Value = 133.963
$ws.Cells.Item(75, 3).Value = -101.8721
$ws.Cells.Item(75, 4).Value = "top"
$ws.Cells.Item(75, 5).Value = -150

$ws.Cells.Item(76, 1).Value = "P6"
$ws.Cells.Item(76, 2).Value = 101.8721
$ws.Cells.Item(76, 3).Value = -133.963
$ws.Cells.Item(76, 4).Value = "top"
$ws.Cells.Item(76, 5).Value = -120

$ws.Cells.Item(77, 1).Value = "P7"
$ws.Cells.Item(77, 2).Value = 90.126
$ws.Cells.Item(77, 3).Value = -177.8
$ws.Cells.Item(77, 4).Value = "top"
$ws.Cells.Item(77, 5).Value = -90

$ws.Cells.Item(78, 1).Value = "P8"
$ws.Cells.Item(78, 2).Value = 101.8721
$ws.Cells.Item(78, 3).Value = -221.637
$ws.Cells.Item(78, 4).Value = "top"
$ws.Cells.Item(78, 5).Value = -60

$ws.Cells.Item(79, 1).Value = "P9"
$ws.Cells.Item(79, 2).Value = 133.963
$ws.Cells.Item(79, 3).Value = -253.7279
$ws.Cells.Item(79, 4).Value = "top"
$ws.Cells.Item(79, 5).Value = -30

$ws.Cells.Item(80, 1).Value = "P10"
$ws.Cells.Item(80, 2).Value = 177.8
$ws.Cells.Item(80, 3).Value = -265.474
$ws.Cells.Item(80, 4).Value = "top"
$ws.Cells.Item(80, 5).Value = 0

$ws.Cells.Item(81, 1).Value = "P11"
$ws.Cells.Item(81, 2).Value = 221.637
$ws.Cells.Item(81, 3).Value = -253.7279
$ws.Cells.Item(81, 4).Value = "top"
$ws.Cells.Item(81, 5).Value = 30

$ws.Cells.Item(82, 1).Value = "P12"
$ws.Cells.Item(82, 2).Value = 253.7279
$ws.Cells.Item(82, 3).Value = -221.637
$ws.Cells.Item(82, 4).Value = "top"
$ws.Cells.Item(82, 5).Value = 60

$ws.Cells.Item(83, 1).Value = "SW1"
$ws.Cells.Item(83, 2).Value = 299.5
$ws.Cells.Item(83, 3).Value = -86.233
$ws.Cells.Item(83, 4).Value = "top"
$ws.Cells.Item(83, 5).Value = 90

# Restore the scrolled/selected view state recorded in the saved workbook.
[void]$ws.Range("H14").Select()
